# RBC birth and death table.xlsx - apply "Checked" column + content updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "Checked" column to the table (4th column) ---
$tbl = $ws.ListObjects.Item(1)
$tbl.ListColumns.Add() | Out-Null

# --- Cell values: paper / birth / death content updates + new Checked column ---
$ws.Range("A1").Value = 'Paper'
$ws.Range("B1").Value = 'Birth'
$ws.Range("C1").Value = 'Death'
$ws.Range("D1").Value = 'Checked'
$ws.Range("A2").Value = 'Antia'
$ws.Range("B2").Value = 'Non-linear, lagged total RBC density'
$ws.Range("C2").Value = 'Constant background mortality'
$ws.Range("A3").Value = 'Greischar'
$ws.Range("B3").Value = 'Linear, based on current E'
$ws.Range("C3").Value = 'Constant background mortality of uninfected RBCs; removal of pRBCs by constant background mortality and saturating immunity'
$ws.Range("A4").Value = 'Haydon'
$ws.Range("B4").Value = 'Linear, based on current E'
$ws.Range("C4").Value = 'NOT SURE WHAT"S GOING ON WITH IMMUNE REMOVAL OF uRBCS; removal of pRBCs by first-order interaction with immune cells (beta*P*I))'
$ws.Range("A5").Value = 'Jakeman'
$ws.Range("B5").Value = 'Considered time varying production of erythrocytes (Zi), constant production of erythrocytes as well as complete dyserythropoiesis (d=0)'
$ws.Range("C5").Value = 'Di is number of uninfected erythrocytes destroyed during period I; they estimate the average ratio of the destruction of uRBCs to pRBCs (d); proportion of pRBCs removed by either immunity or normal destruction of old RBCs; no destruction of uninfected erythrocytes'
$ws.Range("A6").Value = 'Kamiya'
$ws.Range("B6").Value = "Assume baseline replenishment is function of RBC density at homeostatic equilibrium Rc, times survival given background mortality during infection, given by mu_R''. Also have density-dependent replenishment."
$ws.Range("C6").Value = 'Changing rates of general RBC clearance and targeted pRBC clearance'
$ws.Range("A7").Value = 'Khoury'
$ws.Range("B7").Value = 'Linear, either irrespective of E or based on E'
$ws.Range("C7").Value = 'Constant background mortality of uninfected RBCs; merozoites removed at constant rate; no removal of infected RBCs, constant removal or first order interaction with immune cells stimulated by presence of parasitised cells/merozoites'
$ws.Range("D7").Value = 'Yes'
$ws.Range("A8").Value = 'Kochin'
$ws.Range("B8").Value = 'NA'
$ws.Range("C8").Value = 'Mass action term of removal pRBCs via innate immunity kIP'
$ws.Range("D8").Value = 'Yes'
$ws.Range("A9").Value = 'Lim'
$ws.Range("B9").Value = 'Constant input b'
$ws.Range("C9").Value = 'Mass action mu*S removal of susceptible RBC capturing aging; no removal of pRBCs via immunity; include linear mortaliy of uRBCs to capture haemolytic anaemia '
$ws.Range("D9").Value = 'Yes'
$ws.Range("A10").Value = 'McQueen and McKenzie'
$ws.Range("B10").Value = 'RBC source has own ODE, with maximum rate of production; RBC production regulated by rate of loss of RBCs by means of normal sensecence or infection'
$ws.Range("C10").Value = 'No removal of uninfected RBCs via immunity, no removal of infected RBCs via immunity; removal of RBCs via senescence (end of progression through life stages)'
$ws.Range("A11").Value = 'Metcalf'
$ws.Range("A12").Value = 'Mideo'
$ws.Range("B12").Value = 'Linear, based on lagged E'
$ws.Range("C12").Value = 'Constant death rate of RBCs per day "d"'
$ws.Range("D12").Value = 'Yes'
$ws.Range("A13").Value = 'Savill'
$ws.Range("B13").Value = 'PHZ induced anaemia only'
$ws.Range("C13").Value = 'PHZ induced anaemia only'
$ws.Range("A14").Value = 'Wale'
$ws.Range("B14").Value = 'Changing production of reticulocytes'
$ws.Range("C14").Value = 'Changing removal of general RBC clearance and targeted pRBC clearance'
$ws.Range("D14").Value = 'Yes'
$ws.Range("A15").Value = 'Gravenor'
$ws.Range("B15").Value = 'Constant production of erythrocytes'
$ws.Range("C15").Value = 'No destruction of uninfected erythrocytes, just natural death'
$ws.Range("D15").Value = 'Yes'

# --- Wrap text for the whole new Checked column (matches rest of table formatting) ---
$ws.Range("D1:D15").WrapText = $true

# --- Highlight two rows (Haydon / McQueen and McKenzie) in red, matching new review-style font ---
# (Haydon row: only Paper/Birth/Death are highlighted, Checked stays normal;
#  McQueen and McKenzie row: Checked cell is also highlighted, matching source file)
$ws.Range("A4:C4").Font.Color = 255
$ws.Range("A10:D10").Font.Color = 255

# --- Row height adjustments for re-wrapped rows ---
$ws.Rows.Item(6).RowHeight = 68
$ws.Rows.Item(12).RowHeight = 17
$ws.Rows.Item(14).RowHeight = 34
$ws.Rows.Item(15).RowHeight = 34

# --- Column widths for the new Checked column (and spare column E) ---
$ws.Range("D1:E1").ColumnWidth = 17.25

# --- View state: zoom + selection ---
$ws.Activate()
$ws.Range("B5").Select()
$excel.ActiveWindow.Zoom = 156

Write-Output "RBC table updated"
